$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "မကျီးပင်"
$ws.Range("A2").Value = "မန်ကျီးပင်"
$ws.Range("A3").Value = "မန်ကြီးပင်"
$ws.Range("A4").Value = "မက်ကြီးပင်"
$ws.Range("A5").Value = "မန်ကျည်းပင်"
